$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 19.95578266666667
$ws.Range("H2").Value = 59.867348
$ws.Range("I2").Value = 0.0117373419656925
$ws.Range("J2").Value = 0.0117373419656925
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.3116673333333334
$ws.Range("N2").Value = 0.935002
$ws.Range("O2").Value = 0.0414413620607491
$ws.Range("P2").Value = 0.0414413620607491
$ws.Range("Q2").Value = 6.219565568299556
$ws.Range("R2").Value = 55.976090114696
$ws.Range("S2").Value = 0.0004864114380310872
$ws.Range("T2").Value = 0.0004864114380310872
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 19.95578266666667
$ws.Range("H3").Value = 59.867348
$ws.Range("I3").Value = 0.0117373419656925
$ws.Range("J3").Value = 0.0117373419656925
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.794584666666667
$ws.Range("N3").Value = 11.383754
$ws.Range("O3").Value = 0.5045532214096876
$ws.Range("P3").Value = 0.5045532214096876
$ws.Range("Q3").Value = 75.72390691826578
$ws.Range("R3").Value = 681.515162264392
$ws.Range("S3").Value = 0.005922113699577264
$ws.Range("T3").Value = 0.005922113699577264
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 19.95578266666667
$ws.Range("H4").Value = 59.867348
$ws.Range("I4").Value = 0.0117373419656925
$ws.Range("J4").Value = 0.0117373419656925
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.414430666666667
$ws.Range("N4").Value = 10.243292
$ws.Range("O4").Value = 0.4540054165295632
$ws.Range("P4").Value = 0.4540054165295633
$ws.Range("Q4").Value = 68.13763631440177
$ws.Range("R4").Value = 613.238726829616
$ws.Range("S4").Value = 0.005328816828084144
$ws.Range("T4").Value = 0.005328816828084145
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1637.343343333333
$ws.Range("H5").Value = 4912.03003
$ws.Range("I5").Value = 0.9630320723052701
$ws.Range("J5").Value = 0.9630320723052702
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.3116673333333334
$ws.Range("N5").Value = 0.935002
$ws.Range("O5").Value = 0.0414413620607491
$ws.Range("P5").Value = 0.0414413620607491
$ws.Range("Q5").Value = 510.3064335677844
$ws.Range("R5").Value = 4592.75790211006
$ws.Range("S5").Value = 0.0399093607845162
$ws.Range("T5").Value = 0.03990936078451621
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1637.343343333333
$ws.Range("H6").Value = 4912.03003
$ws.Range("I6").Value = 0.9630320723052701
$ws.Range("J6").Value = 0.9630320723052702
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.794584666666667
$ws.Range("N6").Value = 11.383754
$ws.Range("O6").Value = 0.5045532214096876
$ws.Range("P6").Value = 0.5045532214096876
$ws.Range("Q6").Value = 6213.037944681402
$ws.Range("R6").Value = 55917.34150213261
$ws.Range("S6").Value = 0.4859009344024712
$ws.Range("T6").Value = 0.4859009344024713
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1637.343343333333
$ws.Range("H7").Value = 4912.03003
$ws.Range("I7").Value = 0.9630320723052701
$ws.Range("J7").Value = 0.9630320723052702
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.414430666666667
$ws.Range("N7").Value = 10.243292
$ws.Range("O7").Value = 0.4540054165295632
$ws.Range("P7").Value = 0.4540054165295633
$ws.Range("Q7").Value = 5590.595323339861
$ws.Range("R7").Value = 50315.35791005876
$ws.Range("S7").Value = 0.4372217771182826
$ws.Range("T7").Value = 0.4372217771182827
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 17.50081933333334
$ws.Range("H8").Value = 52.502458
$ws.Range("I8").Value = 0.01029341242216722
$ws.Range("J8").Value = 0.01029341242216722
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.3116673333333334
$ws.Range("N8").Value = 0.935002
$ws.Range("O8").Value = 0.0414413620607491
$ws.Range("P8").Value = 0.0414413620607491
$ws.Range("Q8").Value = 5.454433692768445
$ws.Range("R8").Value = 49.089903234916
$ws.Range("S8").Value = 0.000426573031027644
$ws.Range("T8").Value = 0.0004265730310276441
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 17.50081933333334
$ws.Range("H9").Value = 52.502458
$ws.Range("I9").Value = 0.01029341242216722
$ws.Range("J9").Value = 0.01029341242216722
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.794584666666667
$ws.Range("N9").Value = 11.383754
$ws.Range("O9").Value = 0.5045532214096876
$ws.Range("P9").Value = 0.5045532214096876
$ws.Range("Q9").Value = 66.40834069637023
$ws.Range("R9").Value = 597.675066267332
$ws.Range("S9").Value = 0.005193574396902965
$ws.Range("T9").Value = 0.005193574396902966
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 17.50081933333334
$ws.Range("H10").Value = 52.502458
$ws.Range("I10").Value = 0.01029341242216722
$ws.Range("J10").Value = 0.01029341242216722
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.414430666666667
$ws.Range("N10").Value = 10.243292
$ws.Range("O10").Value = 0.4540054165295632
$ws.Range("P10").Value = 0.4540054165295633
$ws.Range("Q10").Value = 59.75533422352623
$ws.Range("R10").Value = 537.7980080117361
$ws.Range("S10").Value = 0.004673264994236608
$ws.Range("T10").Value = 0.00467326499423661
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 25.39612333333333
$ws.Range("H11").Value = 76.18836999999999
$ws.Range("I11").Value = 0.01493717330687017
$ws.Range("J11").Value = 0.01493717330687017
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.3116673333333334
$ws.Range("N11").Value = 0.935002
$ws.Range("O11").Value = 0.0414413620607491
$ws.Range("P11").Value = 0.0414413620607491
$ws.Range("Q11").Value = 7.915142036304444
$ws.Range("R11").Value = 71.23627832673999
$ws.Range("S11").Value = 0.0006190168071741635
$ws.Range("T11").Value = 0.0006190168071741635
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 25.39612333333333
$ws.Range("H12").Value = 76.18836999999999
$ws.Range("I12").Value = 0.01493717330687017
$ws.Range("J12").Value = 0.01493717330687017
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 3.794584666666667
$ws.Range("N12").Value = 11.383754
$ws.Range("O12").Value = 0.5045532214096876
$ws.Range("P12").Value = 0.5045532214096876
$ws.Range("Q12").Value = 96.36774019344222
$ws.Range("R12").Value = 867.3096617409799
$ws.Range("S12").Value = 0.00753659891073614
$ws.Range("T12").Value = 0.00753659891073614
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 25.39612333333333
$ws.Range("H13").Value = 76.18836999999999
$ws.Range("I13").Value = 0.01493717330687017
$ws.Range("J13").Value = 0.01493717330687017
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.414430666666667
$ws.Range("N13").Value = 10.243292
$ws.Range("O13").Value = 0.4540054165295632
$ws.Range("P13").Value = 0.4540054165295633
$ws.Range("Q13").Value = 86.71330232378222
$ws.Range("R13").Value = 780.4197209140399
$ws.Range("S13").Value = 0.006781557588959864
$ws.Range("T13").Value = 0.006781557588959865
